$d = $word.ActiveDocument

# --- Change 1: "(NOMBRE_APP)" -> "TUPARCHE" (title run) ---
$r1 = $d.Content
$found1 = $r1.Find.Execute("(NOMBRE_APP)", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $xml1 = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" + `
        "<w:r><w:rPr><w:rFonts w:cstheme=""minorHAnsi""/><w:b/><w:sz w:val=""24""/><w:szCs w:val=""24""/><w:lang w:val=""es-ES""/></w:rPr><w:t>TUPARCHE</w:t></w:r>" + `
        "</w:p>"
    $r1.InsertXML($xml1)
}

# --- Change 2: split "La capa lógica (BackEnd) debe ser alojada en un servidor de Linux"
#     into three runs, wrapping "BackEnd" with spellcheck proofErr markers ---
$r2 = $d.Content
$found2 = $r2.Find.Execute("La capa lógica (BackEnd) debe ser alojada en un servidor de Linux", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $start2 = $r2.Start
    $r2.Text = ""
    $r3 = $d.Range($start2, $start2)
    $xml2 = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" + `
        "<w:r><w:rPr><w:rFonts w:cstheme=""minorHAnsi""/><w:sz w:val=""24""/><w:szCs w:val=""24""/><w:lang w:val=""es-ES""/></w:rPr><w:t>La capa lógica (</w:t></w:r>" + `
        "<w:proofErr w:type=""spellStart""/>" + `
        "<w:r><w:rPr><w:rFonts w:cstheme=""minorHAnsi""/><w:sz w:val=""24""/><w:szCs w:val=""24""/><w:lang w:val=""es-ES""/></w:rPr><w:t>BackEnd</w:t></w:r>" + `
        "<w:proofErr w:type=""spellEnd""/>" + `
        "<w:r><w:rPr><w:rFonts w:cstheme=""minorHAnsi""/><w:sz w:val=""24""/><w:szCs w:val=""24""/><w:lang w:val=""es-ES""/></w:rPr><w:t>) debe ser alojada en un servidor de Linux</w:t></w:r>" + `
        "</w:p>"
    $r3.InsertXML($xml2)
}
